$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.892.14'
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").Value = '1.878.82'
$ws.Range("E3").Value = '  -1.96%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").Value = '''324.36'
$ws.Range("E5").Value = '  -1.17%  '

$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  -0.36%  '

$ws.Range("D7").Value = '''0.4616'
$ws.Range("E7").Value = '  -1.08%  '

$ws.Range("E8").Value = '  -2.39%  '

$ws.Range("D9").Value = '''0.07856'
$ws.Range("E9").Value = '  -2.22%  '

$ws.Range("D10").Value = '''0.9835'

$ws.Range("E11").Value = '  -2.13%  '

$ws.Range("D12").Value = '1.848.67'
$ws.Range("E12").Value = '  -3.35%  '

$ws.Range("D13").Value = '''6.990'
$ws.Range("E13").Value = '  -2.32%  '

$ws.Range("E14").Value = '  -2.65%  '

$ws.Range("D15").Value = '''0.06981'
$ws.Range("E15").Value = '  +0.35%  '

$ws.Range("D16").Value = '''88.49'
$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("E17").Value = '  -0.34%  '

$ws.Range("D18").Value = '''0.000009939'
$ws.Range("E18").Value = '  -2.46%  '

$ws.Range("D19").Value = '''16.96'
$ws.Range("E19").Value = '  -2.26%  '

$ws.Range("E20").Value = '  -0.36%  '

$ws.Range("D21").Value = '28.893.35'
$ws.Range("E21").Value = '  -1.18%  '

$ws.Range("D22").Value = '''5.263'
$ws.Range("E22").Value = '  -2.81%  '

$ws.Range("E23").Value = '  -1.68%  '

$ws.Range("D24").Value = '''2.103'
$ws.Range("E24").Value = '  +2.17%  '

$ws.Range("D25").Value = '''156.21'
$ws.Range("E25").Value = '  +0.41%  '

$ws.Range("E26").Value = '  -2.07%  '

$ws.Range("D27").Value = '''5.899'
$ws.Range("E27").Value = '  -0.65%  '

$ws.Range("D28").Value = '''117.75'
$ws.Range("E28").Value = '  -2.72%  '

$ws.Range("D29").Value = '''1.904'
$ws.Range("E29").Value = '  -6.90%  '

$ws.Range("D30").Value = '''0.09362'
$ws.Range("E30").Value = '  -0.52%  '

$ws.Range("D31").Value = '''0.9015'
$ws.Range("E31").Value = '  -4.78%  '

$ws.Range("D32").Value = '''5.271'
$ws.Range("E32").Value = '  -1.95%  '

$ws.Range("D33").Value = '''1.318'
$ws.Range("E33").Value = '  -2.76%  '

$ws.Range("E34").Value = '  -0.95%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.05743'
$ws.Range("E35").Value = '  -2.55%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").Value = '''1.170'
$ws.Range("E36").Value = '  -0.49%  '

$ws.Range("D37").Value = '''0.02075'
$ws.Range("E37").Value = '  -1.75%  '

$ws.Range("E38").Value = '  -0.41%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '''7.625'
$ws.Range("E39").Value = '  -6.58%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.5656'
$ws.Range("E40").Value = '  -3.45%  '

$ws.Range("D41").Value = '''0.1774'
$ws.Range("E41").Value = '  -2.69%  '

$ws.Range("D42").Value = '''9.688'
$ws.Range("E42").Value = '  -4.06%  '

$ws.Range("D43").Value = '''11.98'
$ws.Range("E43").Value = '  -0.58%  '

$ws.Range("D44").Value = '''2.234'
$ws.Range("E44").Value = '  -3.72%  '

$ws.Range("D45").Value = '''0.5337'
$ws.Range("E45").Value = '  -2.56%  '

$ws.Range("D46").Value = '''0.07039'
$ws.Range("E46").Value = '  -2.46%  '

$ws.Range("D47").Value = '''1.844'
$ws.Range("E47").Value = '  -4.42%  '

$ws.Range("D48").Value = '''2.541'
$ws.Range("E48").Value = '  +1.56%  '

$ws.Range("D49").Value = '''112.37'
$ws.Range("E49").Value = '  -0.90%  '

$ws.Range("D50").Value = '''1.068'
$ws.Range("E50").Value = '  -5.54%  '

$ws.Range("D51").Value = '''70.83'
$ws.Range("E51").Value = '  -1.73%  '

